$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 through 5 (the old "ECs" related rows), keep only header + row2
$ws.Rows("3:5").Delete()

# Update remaining data row (row 2) with new values
$ws.Range("A2").Value = "Resolving-Mac"
$ws.Range("B2").Value = "Ccl12"
$ws.Range("C2").Value = "Ccr3"
$ws.Range("D2").Value = "Resolving-Mac"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 20.23247666666667
$ws.Range("H2").Value = 60.69743
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1790523333333333
$ws.Range("N2").Value = 0.537157
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 3.622672156278889
$ws.Range("R2").Value = 32.60404940651
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1

$wb.Save()
